# Loan RBI, Variable Instalments
# Insert a new (blank) column before column N on the "Repayment schedule"
# sheet, shifting the existing "Late" / "heading" / "Outstanding" columns
# one slot to the right, and make that sheet the active one (tab + cell
# selection) as left by the author after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Bring this sheet to the front (sets tabSelected / activeTab).
$ws.Activate()

# Insert a new blank column at N, pushing Late/heading/Outstanding -> O/P/Q.
$ws.Columns("N").Insert()

# New column inherits the width of its left neighbour (column M).
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Leave the selection where the author left it.
$ws.Range("R10").Select()
